$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.715.29"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").Value = "1.924.12"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.96"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.698"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.37"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.02"
$ws.Range("E9").Value = "  +10.40%  "
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0996"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.50"
$ws.Range("E13").Value = "  +7.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.801"
$ws.Range("E14").Value = "  +4.74%  "
$ws.Range("D15").Value = "2.205.36"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "1.929.96"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "36.743.53"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.55"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "0.0₃0869"
$ws.Range("E20").Value = "  +5.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.53"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.32"
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("E23").Value = "  +5.41%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.09"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.81"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.84"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  +6.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0611"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0851"
$ws.Range("E36").Value = "  +20.63%  "
$ws.Range("E37").Value = "  -13.54%  "
$ws.Range("B38").Value = "Gas"
$ws.Range("C38").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.11"
$ws.Range("E38").Value = "  +49.26%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.865"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.77"
$ws.Range("E41").Value = "  +12.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0229"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.18"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.11"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "1.343.37"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0812"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("D51").Value = "2.109.11"
$ws.Range("E51").Value = "  +2.21%  "
